# edit.ps1 - rewrite the lab-report body to the "finalized" version:
#  - add a title block (names / course / assignment / date)
#  - rework the two original paragraphs into a fuller write-up
#  - add run instructions + a testing/limitations paragraph
#  - normalize paragraph spacing (after=0, single line) and first-line
#    indents on the "body" paragraphs, matching the target layout.
#
# NOTE: the source also drops the document's header (word/header1.xml /
# the <w:headerReference> in sectPr). There is no supported Word object
# model call that deletes a header/section-header reference outright
# (HeaderFooter.Exists is read-only, and this is the only section, so
# LinkToPrevious has nothing to link to) -- that part of the change is
# left alone.

$d = $word.ActiveDocument

# wdCollapseStart / wdCollapseEnd
$wdCollapseStart = 1
$wdCollapseEnd = 0

# ---------------------------------------------------------------------
# Phase 1: build the paragraph skeleton (13 paragraphs total) while
# everything is still in its default, un-formatted state. Doing the
# structural inserts *before* touching any paragraph formatting avoids
# new paragraphs inheriting indents/spacing from whichever paragraph
# they were split off from.
# ---------------------------------------------------------------------

$pIntro = $d.Paragraphs(1)

# 5 new paragraphs ahead of the intro paragraph:
#   Andrew Staffieri and Benjamin Brouse / CS 281 / Programming
#   Assignment 2 / 11/30/11 / (blank line)
for ($i = 0; $i -lt 5; $i++) {
    $pIntro.Range.InsertParagraphBefore()
}

# Original intro paragraph ("We were tasked...") is now #6.
$pIntro = $d.Paragraphs(6)

# Blank line after the intro paragraph.
$pIntro.Range.InsertParagraphAfter()

# Original second paragraph ("This was done...") is now #8.
$pEnv = $d.Paragraphs(8)

# New paragraph with run instructions, right after it.
$pEnv.Range.InsertParagraphAfter()

# Blank line before the testing/limitations paragraph.
$pRun = $d.Paragraphs(9)
$pRun.Range.InsertParagraphAfter()

# Original (empty) bookmark paragraph is now #11; leave it as-is for now
# (text gets inserted ahead of the bookmark in phase 2).

# Two trailing blank paragraphs at the very end.
$pTest = $d.Paragraphs(11)
$pTest.Range.InsertParagraphAfter()
$pTrail1 = $d.Paragraphs(12)
$pTrail1.Range.InsertParagraphAfter()

Write-Host ('Paragraph count after skeleton build: ' + $d.Paragraphs.Count)

# ---------------------------------------------------------------------
# Phase 2: fill in text + apply final formatting, addressing paragraphs
# strictly by their final (now-stable) index.
# ---------------------------------------------------------------------

function Set-NoExtraSpacing($para) {
    $para.Format.SpaceAfter = 0
    $para.Format.LineSpacingRule = 0
}

function Set-FirstLineIndent($para) {
    $para.Format.FirstLineIndent = 36
}

# 1: Andrew Staffieri and Benjamin Brouse
$p = $d.Paragraphs(1)
$p.Range.Text = 'Andrew Staffieri and Benjamin Brouse'
Set-NoExtraSpacing $p

# 2: CS 281
$p = $d.Paragraphs(2)
$p.Range.Text = 'CS 281'
Set-NoExtraSpacing $p

# 3: Programming Assignment 2
$p = $d.Paragraphs(3)
$p.Range.Text = 'Programming Assignment 2'
Set-NoExtraSpacing $p

# 4: 11/30/11
$p = $d.Paragraphs(4)
$p.Range.Text = '11/30/11'
Set-NoExtraSpacing $p

# 5: (blank)
$p = $d.Paragraphs(5)
Set-NoExtraSpacing $p

# 6: We were tasked ... our alu).
$p = $d.Paragraphs(6)
$p.Range.Text = 'We were tasked (as a group) to create an 8 bit multiplier in VHDL.   The multiplier is the 3-rd version we discussed in class, which has a finite state machine and a clock (except we were given an adder to use instead of our alu).'
Set-NoExtraSpacing $p
Set-FirstLineIndent $p

# 7: (blank, indented)
$p = $d.Paragraphs(7)
Set-NoExtraSpacing $p
Set-FirstLineIndent $p

# 8: The code submitted with this lab ...
$p = $d.Paragraphs(8)
$p.Range.Text = 'The code submitted with this lab was created in Sonata VHDL on a Windows 7 64-bit machine.'
Set-NoExtraSpacing $p
Set-FirstLineIndent $p

# 9: To run the files in this assignment ...
$p = $d.Paragraphs(9)
$p.Range.Text = 'To run the files in this assignment you must have access to a VHDL simulator (i.e. Sonata).  Load the appropriate files in the simulator (all .vhd files included with this submission) and add the appropriate level to the top-level (multiplier tester).  You can now simulate using the test bench included.'
Set-NoExtraSpacing $p

# 10: (blank, indented)
$p = $d.Paragraphs(10)
Set-NoExtraSpacing $p
Set-FirstLineIndent $p

# 11: Testing for this assignment ... (bookmark paragraph; insert the
# text ahead of the existing _GoBack bookmark so it keeps trailing the
# final run, exactly as in the source document).
$p = $d.Paragraphs(11)
$r = $p.Range
$r.Collapse($wdCollapseStart)
$r.InsertBefore('Testing for this assignment was a frustrating a painstaking task for us.  It was tough to time when one multiplication process would be done and when another would begin (so we could do a reset on the HILO registers).  Our code works perfectly except for one flaw: when the least significant bit is 1 and the mcand should be added to the HI register, our code shifts first and then adds.  We worked for an extremely long time attempting to overcome this problem but could not.  This is the only reason our output values are just a little off.')
Set-NoExtraSpacing $p
Set-FirstLineIndent $p

# 12 & 13: trailing blanks (no indent)
$p = $d.Paragraphs(12)
Set-NoExtraSpacing $p

$p = $d.Paragraphs(13)
Set-NoExtraSpacing $p

Write-Host ('Final paragraph count: ' + $d.Paragraphs.Count)
